$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.361.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.47%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.796.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.47%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.78%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.27%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3776"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3460"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.51%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.37"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.43%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.202"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.97%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07526"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.27%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.006"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.25%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.76"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.86%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.470"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.70%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.803.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.62%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.037"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.42%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001095"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.68%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06671"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.95%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "84.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.10%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.006"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.39%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.440"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.78%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.436.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.91%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.88%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.443"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.23%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.578"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.52%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +10.62%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.440"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.68%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "149.94"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.28%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.010.54"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.87%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "133.66"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.50%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.075"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.073"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.48%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08677"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.33%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.29"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.49%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.675"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.41%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.447"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.54%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6859"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +11.59%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.923"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.89%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02354"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.91%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2203"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.67%  "

$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.06345"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.70%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.286"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.35%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.01%  "

$ws.Range("E45").Value = "  +0.37%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6380"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.63%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.843"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.42%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.132"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.42%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "131.02"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.87%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07222"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.83%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.91%  "
